$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the AveragePeriod value (column N, row 2) from 1000 to 200
$ws.Range("N2").Value = 200

# Update the active selection to N3 (matches the recorded cursor position in the file)
$ws.Range("N3").Select()
